# Applies the "Results from R script" update to the PAL.MI sheet:
#  - fixes the timestamp and low price recorded for the existing last row (76)
#  - appends a brand-new row (77) with the next day's OHLC/volume data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 76: correct the date serial and the low ("D") value -----------------
$ws.Range("A76").Value = 45447.2916666667
$ws.Range("D76").Value = 6.42000007629395

# --- Row 77: brand-new observation -------------------------------------------
$ws.Range("A77").Value = 45448.6229513889
$ws.Range("B77").Value = 2100
$ws.Range("C77").Value = 6.40000009536743
$ws.Range("D77").Value = 6.30000019073486
$ws.Range("E77").Value = 6.40000009536743
$ws.Range("F77").Value = 6.30000019073486

# adj_close ("G") is stored as text in this sheet (shared string), same as
# every other row, so force a text format before writing the numeric-looking
# string, then copy row 76's (General) formatting back on top so the cell
# keeps the same default styling as the rest of the column.
$ws.Range("G77").NumberFormat = "@"
$ws.Range("G77").Value = "6.30000019073486"
$ws.Range("G76").Copy()
$ws.Range("G77").PasteSpecial(-4122)

# ticker ("H") is always the same text value already used on every other row;
# copy it straight from the row above so it reuses the same shared string
# (Value2 is used for reading because Value as a getter is unreliable here).
$ws.Range("H77").Value = $ws.Range("H76").Value2

# Match the date column's custom date/time number format (and underlying
# cell style) used throughout column A by copying it from the row above.
$ws.Range("A76").Copy()
$ws.Range("A77").PasteSpecial(-4122)

$excel.CutCopyMode = 0
